$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new LEM/LED rows below the existing data (rows 83-86),
# entering the A column first and then the B column, matching the
# order in which the new unique strings were introduced.
$ws.Range("A83").Value = "LEM-275-16-2722KS"
$ws.Range("A84").Value = "LEM-293-00-30KS"
$ws.Range("A85").Value = "LEM-307-00-35KH"
$ws.Range("A86").Value = "LEM-319-00-27KU"

$ws.Range("B83").Value = "LED-275-S35-2722"
$ws.Range("B84").Value = "LED-293-S00-30"
$ws.Range("B85").Value = "LED-307-H00-35"
$ws.Range("B86").Value = "LED-319-U00-27"

# Match the centered alignment style used by the rest of column B by
# copying the format from the last existing row instead of setting
# alignment properties directly (which would create a stray style).
$ws.Range("B82").Copy()
$ws.Range("B83:B86").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to reflect where the user ended up after entry
[void]$ws.Range("B87").Select()
